$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53,8).Value = 297.8889  # ALC!H53 (G=5479)
$ws.Cells.Item(53,10).Value = 278.4  # ALC!J53 (G=5479)
$ws.Cells.Item(53,12).Value = 278.4  # ALC!L53 (G=5479)
$ws.Cells.Item(53,14).Value = -1552.4  # ALC!N53 (G=5479)

$ws.Cells.Item(111,8).Value = 1217.2858  # ALC!H111 (G=27768)
$ws.Cells.Item(111,9).Value = 1217.2858  # ALC!I111 (G=27768)
$ws.Cells.Item(111,10).Value = 0  # ALC!J111 (G=27768)
$ws.Cells.Item(111,11).Value = 3651.8574  # ALC!K111 (G=27768)
$ws.Cells.Item(111,12).Value = 0  # ALC!L111 (G=27768)
$ws.Cells.Item(111,13).Value = -584.8574000000003  # ALC!M111 (G=27768)
$ws.Cells.Item(111,14).Value = $null  # ALC!N111 (G=27768)

$ws.Cells.Item(131,8).Value = 1586.1875  # ALC!H131 (G=36108)
$ws.Cells.Item(131,10).Value = 4933.3335  # ALC!J131 (G=36108)
$ws.Cells.Item(131,12).Value = 14800.0005  # ALC!L131 (G=36108)
$ws.Cells.Item(131,14).Value = -24880.0005  # ALC!N131 (G=36108)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4,8).Value = 713.3333  # ARM!H4 (G=5071)
$ws.Cells.Item(4,9).Value = 302.8  # ARM!I4 (G=5071)
$ws.Cells.Item(4,10).Value = 1226.5  # ARM!J4 (G=5071)
$ws.Cells.Item(4,11).Value = 302.8  # ARM!K4 (G=5071)
$ws.Cells.Item(4,12).Value = 1226.5  # ARM!L4 (G=5071)
$ws.Cells.Item(4,13).Value = -186.8  # ARM!M4 (G=5071)
$ws.Cells.Item(4,14).Value = -1458.5  # ARM!N4 (G=5071)

$ws.Cells.Item(6,8).Value = 1750  # ARM!H6 (G=2226)
$ws.Cells.Item(6,9).Value = 1750  # ARM!I6 (G=2226)
$ws.Cells.Item(6,11).Value = 1750  # ARM!K6 (G=2226)
$ws.Cells.Item(6,13).Value = -1577  # ARM!M6 (G=2226)

$ws.Cells.Item(19,8).Value = 2286.5715  # ARM!H19 (G=3550)
$ws.Cells.Item(19,9).Value = 1601.4  # ARM!I19 (G=3550)
$ws.Cells.Item(19,10).Value = 3999.5  # ARM!J19 (G=3550)
$ws.Cells.Item(19,11).Value = 1601.4  # ARM!K19 (G=3550)
$ws.Cells.Item(19,12).Value = 3999.5  # ARM!L19 (G=3550)
$ws.Cells.Item(19,13).Value = -1372.4  # ARM!M19 (G=3550)
$ws.Cells.Item(19,14).Value = -4457.5  # ARM!N19 (G=3550)

$ws.Cells.Item(25,8).Value = 0  # ARM!H25 (G=2471)
$ws.Cells.Item(25,9).Value = 0  # ARM!I25 (G=2471)
$ws.Cells.Item(25,11).Value = 0  # ARM!K25 (G=2471)
$ws.Cells.Item(25,13).Value = $null  # ARM!M25 (G=2471)

$ws.Cells.Item(30,8).Value = 1344.75  # ARM!H30 (G=2712)
$ws.Cells.Item(30,10).Value = 1232.5  # ARM!J30 (G=2712)
$ws.Cells.Item(30,12).Value = 1232.5  # ARM!L30 (G=2712)
$ws.Cells.Item(30,14).Value = -1532.5  # ARM!N30 (G=2712)

$ws.Cells.Item(45,8).Value = 5470.3335  # ARM!H45 (G=27714)
$ws.Cells.Item(45,9).Value = 4012  # ARM!I45 (G=27714)
$ws.Cells.Item(45,11).Value = 4012  # ARM!K45 (G=27714)
$ws.Cells.Item(45,13).Value = -3635  # ARM!M45 (G=27714)

$ws.Cells.Item(61,8).Value = 3006.2856  # ARM!H61 (G=43999)
$ws.Cells.Item(61,9).Value = 2781.6365  # ARM!I61 (G=43999)
$ws.Cells.Item(61,11).Value = 2781.6365  # ARM!K61 (G=43999)
$ws.Cells.Item(61,13).Value = -2569.6365  # ARM!M61 (G=43999)

$ws.Cells.Item(102,8).Value = 4446.8945  # ARM!H102 (G=19945)
$ws.Cells.Item(102,9).Value = 4427.5884  # ARM!I102 (G=19945)
$ws.Cells.Item(102,10).Value = 4611  # ARM!J102 (G=19945)
$ws.Cells.Item(102,11).Value = 4427.5884  # ARM!K102 (G=19945)
$ws.Cells.Item(102,12).Value = 4611  # ARM!L102 (G=19945)
$ws.Cells.Item(102,13).Value = -2805.5884  # ARM!M102 (G=19945)
$ws.Cells.Item(102,14).Value = -7855  # ARM!N102 (G=19945)

$ws.Cells.Item(132,8).Value = 1668515.4  # ARM!H132 (G=43997)
$ws.Cells.Item(132,9).Value = 2085073.8  # ARM!I132 (G=43997)
$ws.Cells.Item(132,10).Value = 2281.6667  # ARM!J132 (G=43997)
$ws.Cells.Item(132,11).Value = 6255221.4  # ARM!K132 (G=43997)
$ws.Cells.Item(132,12).Value = 6845.000100000001  # ARM!L132 (G=43997)
$ws.Cells.Item(132,13).Value = -6252691.4  # ARM!M132 (G=43997)
$ws.Cells.Item(132,14).Value = -11905.0001  # ARM!N132 (G=43997)

$ws.Cells.Item(133,8).Value = 56652.332  # ARM!H133 (G=41857)
$ws.Cells.Item(133,10).Value = 56652.332  # ARM!J133 (G=41857)
$ws.Cells.Item(133,12).Value = 56652.332  # ARM!L133 (G=41857)
$ws.Cells.Item(133,14).Value = -61712.332  # ARM!N133 (G=41857)

$ws.Cells.Item(136,8).Value = 3006.2856  # ARM!H136 (G=43999)
$ws.Cells.Item(136,9).Value = 2781.6365  # ARM!I136 (G=43999)
$ws.Cells.Item(136,11).Value = 8344.9095  # ARM!K136 (G=43999)
$ws.Cells.Item(136,13).Value = -5794.9095  # ARM!M136 (G=43999)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7,8).Value = 546533.5600000001  # BSM!H7 (G=1602)
$ws.Cells.Item(7,9).Value = 667629.9  # BSM!I7 (G=1602)
$ws.Cells.Item(7,11).Value = 667629.9  # BSM!K7 (G=1602)
$ws.Cells.Item(7,13).Value = -667516.9  # BSM!M7 (G=1602)

$ws.Cells.Item(35,8).Value = 30642.572  # BSM!H35 (G=2350)
$ws.Cells.Item(35,10).Value = 30642.572  # BSM!J35 (G=2350)
$ws.Cells.Item(35,12).Value = 30642.572  # BSM!L35 (G=2350)
$ws.Cells.Item(35,14).Value = -31262.572  # BSM!N35 (G=2350)

$ws.Cells.Item(105,8).Value = 2984.25  # BSM!H105 (G=19947)
$ws.Cells.Item(105,9).Value = 2681.682  # BSM!I105 (G=19947)
$ws.Cells.Item(105,11).Value = 2681.682  # BSM!K105 (G=19947)
$ws.Cells.Item(105,13).Value = -934.6819999999998  # BSM!M105 (G=19947)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2,8).Value = 302.6  # CRP!H2 (G=1820)
$ws.Cells.Item(2,9).Value = 4  # CRP!I2 (G=1820)
$ws.Cells.Item(2,10).Value = 501.66666  # CRP!J2 (G=1820)
$ws.Cells.Item(2,11).Value = 4  # CRP!K2 (G=1820)
$ws.Cells.Item(2,12).Value = 501.66666  # CRP!L2 (G=1820)
$ws.Cells.Item(2,13).Value = 109  # CRP!M2 (G=1820)
$ws.Cells.Item(2,14).Value = -727.66666  # CRP!N2 (G=1820)

$ws.Cells.Item(6,8).Value = 354546.1  # CRP!H6 (G=2219)
$ws.Cells.Item(6,10).Value = 125001.5  # CRP!J6 (G=2219)
$ws.Cells.Item(6,12).Value = 125001.5  # CRP!L6 (G=2219)
$ws.Cells.Item(6,14).Value = -125227.5  # CRP!N6 (G=2219)

$ws.Cells.Item(7,8).Value = 173.65  # CRP!H7 (G=5361)
$ws.Cells.Item(7,9).Value = 132.90909  # CRP!I7 (G=5361)
$ws.Cells.Item(7,10).Value = 223.44444  # CRP!J7 (G=5361)
$ws.Cells.Item(7,11).Value = 132.90909  # CRP!K7 (G=5361)
$ws.Cells.Item(7,12).Value = 223.44444  # CRP!L7 (G=5361)
$ws.Cells.Item(7,13).Value = -19.90908999999999  # CRP!M7 (G=5361)
$ws.Cells.Item(7,14).Value = -449.44444  # CRP!N7 (G=5361)

$ws.Cells.Item(10,8).Value = 2073.1428  # CRP!H10 (G=1997)
$ws.Cells.Item(10,9).Value = 1168.6666  # CRP!I10 (G=1997)
$ws.Cells.Item(10,11).Value = 1168.6666  # CRP!K10 (G=1997)
$ws.Cells.Item(10,13).Value = -1029.6666  # CRP!M10 (G=1997)

$ws.Cells.Item(19,8).Value = 389.4  # CRP!H19 (G=2233)
$ws.Cells.Item(19,9).Value = 410.07144  # CRP!I19 (G=2233)
$ws.Cells.Item(19,10).Value = 100  # CRP!J19 (G=2233)
$ws.Cells.Item(19,11).Value = 410.07144  # CRP!K19 (G=2233)
$ws.Cells.Item(19,12).Value = 100  # CRP!L19 (G=2233)
$ws.Cells.Item(19,13).Value = -240.07144  # CRP!M19 (G=2233)
$ws.Cells.Item(19,14).Value = -440  # CRP!N19 (G=2233)

$ws.Cells.Item(22,8).Value = 1608.8695  # CRP!H22 (G=5367)
$ws.Cells.Item(22,9).Value = 567.8333  # CRP!I22 (G=5367)
$ws.Cells.Item(22,10).Value = 2744.5454  # CRP!J22 (G=5367)
$ws.Cells.Item(22,11).Value = 567.8333  # CRP!K22 (G=5367)
$ws.Cells.Item(22,12).Value = 2744.5454  # CRP!L22 (G=5367)
$ws.Cells.Item(22,13).Value = -217.8333  # CRP!M22 (G=5367)
$ws.Cells.Item(22,14).Value = -3444.5454  # CRP!N22 (G=5367)

$ws.Cells.Item(24,8).Value = 389.4  # CRP!H24 (G=2233)
$ws.Cells.Item(24,9).Value = 410.07144  # CRP!I24 (G=2233)
$ws.Cells.Item(24,10).Value = 100  # CRP!J24 (G=2233)
$ws.Cells.Item(24,11).Value = 410.07144  # CRP!K24 (G=2233)
$ws.Cells.Item(24,12).Value = 100  # CRP!L24 (G=2233)
$ws.Cells.Item(24,13).Value = -240.07144  # CRP!M24 (G=2233)
$ws.Cells.Item(24,14).Value = -440  # CRP!N24 (G=2233)

$ws.Cells.Item(31,8).Value = 3439.7693  # CRP!H31 (G=44023)
$ws.Cells.Item(31,9).Value = 3421.7  # CRP!I31 (G=44023)
$ws.Cells.Item(31,11).Value = 3421.7  # CRP!K31 (G=44023)
$ws.Cells.Item(31,13).Value = -3126.7  # CRP!M31 (G=44023)

$ws.Cells.Item(34,8).Value = 3439.7693  # CRP!H34 (G=44023)
$ws.Cells.Item(34,9).Value = 3421.7  # CRP!I34 (G=44023)
$ws.Cells.Item(34,11).Value = 3421.7  # CRP!K34 (G=44023)
$ws.Cells.Item(34,13).Value = -3219.7  # CRP!M34 (G=44023)

$ws.Cells.Item(99,8).Value = 6203262  # CRP!H99 (G=36198)
$ws.Cells.Item(99,9).Value = 21706.273  # CRP!I99 (G=36198)
$ws.Cells.Item(99,10).Value = 19802686  # CRP!J99 (G=36198)
$ws.Cells.Item(99,11).Value = 21706.273  # CRP!K99 (G=36198)
$ws.Cells.Item(99,12).Value = 19802686  # CRP!L99 (G=36198)
$ws.Cells.Item(99,13).Value = -20208.273  # CRP!M99 (G=36198)
$ws.Cells.Item(99,14).Value = -19805682  # CRP!N99 (G=36198)

$ws.Cells.Item(107,8).Value = 554.1177  # CRP!H107 (G=27689)
$ws.Cells.Item(107,9).Value = 494.6154  # CRP!I107 (G=27689)
$ws.Cells.Item(107,11).Value = 494.6154  # CRP!K107 (G=27689)
$ws.Cells.Item(107,13).Value = 1425.3846  # CRP!M107 (G=27689)

$ws.Cells.Item(126,8).Value = 6203262  # CRP!H126 (G=36198)
$ws.Cells.Item(126,9).Value = 21706.273  # CRP!I126 (G=36198)
$ws.Cells.Item(126,10).Value = 19802686  # CRP!J126 (G=36198)
$ws.Cells.Item(126,11).Value = 65118.819  # CRP!K126 (G=36198)
$ws.Cells.Item(126,12).Value = 59408058  # CRP!L126 (G=36198)
$ws.Cells.Item(126,13).Value = -62648.819  # CRP!M126 (G=36198)
$ws.Cells.Item(126,14).Value = -59412998  # CRP!N126 (G=36198)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11,8).Value = 52631800  # CUL!H11 (G=4745)
$ws.Cells.Item(11,9).Value = 148.16667  # CUL!I11 (G=4745)
$ws.Cells.Item(11,11).Value = 444.50001  # CUL!K11 (G=4745)
$ws.Cells.Item(11,13).Value = -304.50001  # CUL!M11 (G=4745)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14,8).Value = 167504460  # GSM!H14 (G=4198)
$ws.Cells.Item(14,9).Value = 201005150  # GSM!I14 (G=4198)
$ws.Cells.Item(14,10).Value = 1000  # GSM!J14 (G=4198)
$ws.Cells.Item(14,11).Value = 201005150  # GSM!K14 (G=4198)
$ws.Cells.Item(14,12).Value = 1000  # GSM!L14 (G=4198)
$ws.Cells.Item(14,13).Value = -201004982  # GSM!M14 (G=4198)
$ws.Cells.Item(14,14).Value = -1336  # GSM!N14 (G=4198)

$ws.Cells.Item(20,8).Value = 21000  # GSM!H20 (G=4095)
$ws.Cells.Item(20,9).Value = 21000  # GSM!I20 (G=4095)
$ws.Cells.Item(20,10).Value = 0  # GSM!J20 (G=4095)
$ws.Cells.Item(20,11).Value = 21000  # GSM!K20 (G=4095)
$ws.Cells.Item(20,12).Value = 0  # GSM!L20 (G=4095)
$ws.Cells.Item(20,13).Value = -20755  # GSM!M20 (G=4095)
$ws.Cells.Item(20,14).Value = $null  # GSM!N20 (G=4095)

$ws.Cells.Item(97,8).Value = 1281.975  # GSM!H97 (G=19940)
$ws.Cells.Item(97,9).Value = 1169.7297  # GSM!I97 (G=19940)
$ws.Cells.Item(97,11).Value = 1169.7297  # GSM!K97 (G=19940)
$ws.Cells.Item(97,13).Value = -673.7297000000001  # GSM!M97 (G=19940)

$ws.Cells.Item(102,8).Value = 1514.2354  # GSM!H102 (G=36169)
$ws.Cells.Item(102,9).Value = 1514.2354  # GSM!I102 (G=36169)
$ws.Cells.Item(102,10).Value = 0  # GSM!J102 (G=36169)
$ws.Cells.Item(102,11).Value = 1514.2354  # GSM!K102 (G=36169)
$ws.Cells.Item(102,12).Value = 0  # GSM!L102 (G=36169)
$ws.Cells.Item(102,13).Value = 107.7646  # GSM!M102 (G=36169)
$ws.Cells.Item(102,14).Value = $null  # GSM!N102 (G=36169)

$ws.Cells.Item(126,8).Value = 6750.2354  # GSM!H126 (G=36184)
$ws.Cells.Item(126,9).Value = 7846.4614  # GSM!I126 (G=36184)
$ws.Cells.Item(126,10).Value = 3187.5  # GSM!J126 (G=36184)
$ws.Cells.Item(126,11).Value = 23539.3842  # GSM!K126 (G=36184)
$ws.Cells.Item(126,12).Value = 9562.5  # GSM!L126 (G=36184)
$ws.Cells.Item(126,13).Value = -21069.3842  # GSM!M126 (G=36184)
$ws.Cells.Item(126,14).Value = -14502.5  # GSM!N126 (G=36184)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22,8).Value = 2482.2964  # LTW!H22 (G=5277)
$ws.Cells.Item(22,9).Value = 934.8570999999999  # LTW!I22 (G=5277)
$ws.Cells.Item(22,10).Value = 3023.9  # LTW!J22 (G=5277)
$ws.Cells.Item(22,11).Value = 934.8570999999999  # LTW!K22 (G=5277)
$ws.Cells.Item(22,12).Value = 3023.9  # LTW!L22 (G=5277)
$ws.Cells.Item(22,13).Value = -639.8570999999999  # LTW!M22 (G=5277)
$ws.Cells.Item(22,14).Value = -3613.9  # LTW!N22 (G=5277)

$ws.Cells.Item(27,8).Value = 2482.2964  # LTW!H27 (G=5277)
$ws.Cells.Item(27,9).Value = 934.8570999999999  # LTW!I27 (G=5277)
$ws.Cells.Item(27,10).Value = 3023.9  # LTW!J27 (G=5277)
$ws.Cells.Item(27,11).Value = 934.8570999999999  # LTW!K27 (G=5277)
$ws.Cells.Item(27,12).Value = 3023.9  # LTW!L27 (G=5277)
$ws.Cells.Item(27,13).Value = -827.8570999999999  # LTW!M27 (G=5277)
$ws.Cells.Item(27,14).Value = -3237.9  # LTW!N27 (G=5277)

$ws.Cells.Item(55,8).Value = 1113.359  # LTW!H55 (G=5284)
$ws.Cells.Item(55,9).Value = 1149.4706  # LTW!I55 (G=5284)
$ws.Cells.Item(55,11).Value = 1149.4706  # LTW!K55 (G=5284)
$ws.Cells.Item(55,13).Value = -976.4706000000001  # LTW!M55 (G=5284)

$ws.Cells.Item(133,8).Value = 60001  # LTW!H133 (G=41903)
$ws.Cells.Item(133,10).Value = 60001  # LTW!J133 (G=41903)
$ws.Cells.Item(133,12).Value = 60001  # LTW!L133 (G=41903)
$ws.Cells.Item(133,14).Value = -65061  # LTW!N133 (G=41903)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81,8).Value = 3276.05  # WVR!H81 (G=12596)
$ws.Cells.Item(81,9).Value = 2068.0908  # WVR!I81 (G=12596)
$ws.Cells.Item(81,10).Value = 4752.4443  # WVR!J81 (G=12596)
$ws.Cells.Item(81,11).Value = 4136.1816  # WVR!K81 (G=12596)
$ws.Cells.Item(81,12).Value = 9504.8886  # WVR!L81 (G=12596)
$ws.Cells.Item(81,13).Value = -3075.1816  # WVR!M81 (G=12596)
$ws.Cells.Item(81,14).Value = -11626.8886  # WVR!N81 (G=12596)

$ws.Cells.Item(84,8).Value = 3276.05  # WVR!H84 (G=12596)
$ws.Cells.Item(84,9).Value = 2068.0908  # WVR!I84 (G=12596)
$ws.Cells.Item(84,10).Value = 4752.4443  # WVR!J84 (G=12596)
$ws.Cells.Item(84,11).Value = 20680.908  # WVR!K84 (G=12596)
$ws.Cells.Item(84,12).Value = 47524.443  # WVR!L84 (G=12596)
$ws.Cells.Item(84,13).Value = -15376.908  # WVR!M84 (G=12596)
$ws.Cells.Item(84,14).Value = -58132.443  # WVR!N84 (G=12596)

$ws.Cells.Item(126,8).Value = 2220.6667  # WVR!H126 (G=36210)
$ws.Cells.Item(126,9).Value = 2090.7778  # WVR!I126 (G=36210)
$ws.Cells.Item(126,10).Value = 3000  # WVR!J126 (G=36210)
$ws.Cells.Item(126,11).Value = 6272.3334  # WVR!K126 (G=36210)
$ws.Cells.Item(126,12).Value = 9000  # WVR!L126 (G=36210)
$ws.Cells.Item(126,13).Value = -3802.3334  # WVR!M126 (G=36210)
